$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 546.4480874316939
$ws.Range("D2").Value = 124.024199843872
$ws.Range("C3").Value = 2371.194379391101
$ws.Range("C4").Value = 0
